$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A34").Value = 12
$ws.Range("B34").Value = "Verify the text"
$ws.Range("C34").Value = "Go to homepage of partyguard"
$ws.Range("D34").Value = "Text should be present "
$ws.Range("E34").Value = "Text is present"
$ws.Range("F34").Value = "Pass"

$ws.Range("C35").Value = "Click on contactus page"

$ws.Range("C36").Value = "verify the text is present"

$ws.Range("A37").Value = 13
$ws.Range("B37").Value = "Login button"
$ws.Range("C37").Value = "Go to homepage of partguard"
$ws.Range("D37").Value = "Redirect to login page"
$ws.Range("E37").Value = "login page is displayed"
$ws.Range("F37").Value = "Pass"

$ws.Range("C38").Value = "Go to contactus page of partyguard"

$ws.Range("C39").Value = "Click on login button"

$ws.Range("A40").Value = 14
$ws.Range("B40").Value = "Check the text fields"
$ws.Range("D40").Value = "Should be able to type"
$ws.Range("E40").Value = "User can be able to type"
$ws.Range("C40").Value = "click on the text field "
$ws.Range("F40").Value = "Pass"

$ws.Range("C41").Value = "text field should be enabled"

$ws.Range("A42").Value = 15
$ws.Range("B42").Value = "Verify the text"
$ws.Range("C42").Value = "Go to contactus page of partyguard"
$ws.Range("D42").Value = "Text should be present "
$ws.Range("E42").Value = "Text is present"
$ws.Range("F42").Value = "Pass"

$ws.Range("C43").Value = "Check for the text present in contactus"

$ws.Range("A44").Value = 16
$ws.Range("B44").Value = "Contact details"
$ws.Range("C44").Value = "Go to home page of partyguard"
$ws.Range("D44").Value = "contact details should be displayed"
$ws.Range("E44").Value = "contact details are displayed"
$ws.Range("F44").Value = "pass"

$ws.Range("C45").Value = "Click on contactus page"

$ws.Range("A46").Value = 17
$ws.Range("B46").Value = "Send message button"
$ws.Range("C46").Value = "Go to home page of partyguard"

$ws.Range("C47").Value = "Click on contactus page"

$ws.Range("C48").Value = "enter the details "

$ws.Range("C49").Value = "click the send message button"

$ws.Range("D46").Value = "message should be sent to email"
$ws.Range("E46").Value = "message is sent"
$ws.Range("F46").Value = "pass"

$ws.Range("A50").Value = 18
$ws.Range("B50").Value = "Empty text field validation"
$ws.Range("C50").Value = "Go to homepage of partyguard"
$ws.Range("D50").Value = "validation should be give to fill the textfield"
$ws.Range("E50").Value = "Validation is given"
$ws.Range("F50").Value = "Pass"

$ws.Range("C51").Value = "click on submit button"

$ws.Range("A52").Value = 19
$ws.Range("B52").Value = "Email functionality"
$ws.Range("C52").Value = "Go to contactus page of partyguard"
$ws.Range("D52").Value = "quote should be sent to the user through email"
$ws.Range("E52").Value = "quote is received through email"
$ws.Range("F52").Value = "Pass"

$ws.Range("C53").Value = "enter your name in the text field"

$ws.Range("C54").Value = "enter the email address "

$ws.Range("C55").Value = "click on submit button"

$ws.Range("A56").Value = 20
$ws.Range("B56").Value = "email text content"
$ws.Range("C56").Value = "Go to homepage of partyguard"
$ws.Range("E56").Value = "Text is viewed"
$ws.Range("F56").Value = "Pass"

$ws.Range("C57").Value = "enter your name in the text field"

$ws.Range("C58").Value = "enter the email address "

$ws.Range("C59").Value = "click on submit button"

$ws.Range("C60").Value = "Check the inbox of your email"

$ws.Range("C61").Value = "view the quote from partyguard service"

$ws.Range("A62").Value = 21
$ws.Range("B62").Value = "Text allignment"
$ws.Range("C62").Value = "Go to homepage of partyguard"
$ws.Range("D62").Value = "Text should be readable and alligned properly"
$ws.Range("E62").Value = "text is alligned properly and readable"
$ws.Range("F62").Value = "Pass"

$ws.Range("C63").Value = "enter your name in the text field"

$ws.Range("C64").Value = "enter the email address "

$ws.Range("C65").Value = "click on submit button"

$ws.Range("C66").Value = "Check the inbox of your email"

$ws.Range("C67").Value = "view the quote from partyguard service"

$ws.Range("C68").Value = "Verify the text allignment"

$ws.Range("D56").Value = "Message text should be viewed by the user"

$excel.ActiveWindow.ScrollRow = 47
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C60").Select()
